$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "29.389.12"
$ws.Range("E2").Value = "  -0.61%  "
Set-TextValue $ws.Range("D3") "1.895.66"
$ws.Range("E3").Value = "  -1.17%  "
Set-TextValue $ws.Range("D4") "1.002"
$ws.Range("E4").Value = "  -0.10%  "
Set-TextValue $ws.Range("D5") "324.08"
$ws.Range("E5").Value = "  -3.29%  "
$ws.Range("E6").Value = "  -0.10%  "
Set-TextValue $ws.Range("D7") "0.4771"
$ws.Range("E7").Value = "  +2.20%  "
Set-TextValue $ws.Range("D8") "0.4056"
$ws.Range("E8").Value = "  -1.75%  "
Set-TextValue $ws.Range("D9") "0.08026"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("E10").Value = "  -1.19%  "
Set-TextValue $ws.Range("D11") "23.34"
$ws.Range("E11").Value = "  +3.92%  "
Set-TextValue $ws.Range("D12") "1.938.26"
$ws.Range("E12").Value = "  -1.77%  "
Set-TextValue $ws.Range("D13") "5.933"
$ws.Range("E13").Value = "  -1.33%  "
Set-TextValue $ws.Range("D14") "7.052"
$ws.Range("E14").Value = "  -1.95%  "
Set-TextValue $ws.Range("D15") "89.55"
$ws.Range("E15").Value = "  -0.52%  "
Set-TextValue $ws.Range("D16") "1.003"
$ws.Range("E16").Value = "  -0.11%  "
Set-TextValue $ws.Range("D17") "0.06677"
$ws.Range("E17").Value = "  +1.19%  "
Set-TextValue $ws.Range("D18") "0.00001024"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("E19").Value = "  -1.23%  "
Set-TextValue $ws.Range("D20") "1.000"
$ws.Range("E20").Value = "  -0.23%  "
Set-TextValue $ws.Range("D21") "29.401.50"
$ws.Range("E21").Value = "  -0.53%  "
Set-TextValue $ws.Range("D22") "5.521"
Set-TextValue $ws.Range("D23") "11.69"
$ws.Range("E23").Value = "  +0.44%  "
Set-TextValue $ws.Range("D24") "2.154"
$ws.Range("E24").Value = "  -2.16%  "
Set-TextValue $ws.Range("D25") "2.109.92"
$ws.Range("E25").Value = "  -4.14%  "
Set-TextValue $ws.Range("D26") "154.37"
$ws.Range("E26").Value = "  -1.64%  "
Set-TextValue $ws.Range("D27") "19.70"
$ws.Range("E27").Value = "  -1.14%  "
Set-TextValue $ws.Range("D28") "6.050"
$ws.Range("E28").Value = "  +5.39%  "
Set-TextValue $ws.Range("D29") "2.085"
$ws.Range("E29").Value = "  -2.76%  "
Set-TextValue $ws.Range("D30") "117.95"
$ws.Range("E30").Value = "  +0.33%  "
Set-TextValue $ws.Range("D31") "1.019"
$ws.Range("E31").Value = "  -4.91%  "
Set-TextValue $ws.Range("D32") "0.09496"
$ws.Range("E32").Value = "  +0.15%  "
Set-TextValue $ws.Range("D33") "1.384"
$ws.Range("E33").Value = "  -3.50%  "
Set-TextValue $ws.Range("D34") "3.525"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("E35").Value = "  -1.28%  "
Set-TextValue $ws.Range("D36") "0.02244"
$ws.Range("E36").Value = "  -1.12%  "
Set-TextValue $ws.Range("D37") "0.06043"
$ws.Range("E37").Value = "  -1.63%  "
Set-TextValue $ws.Range("D38") "1.169"
$ws.Range("E38").Value = "  -1.04%  "
Set-TextValue $ws.Range("D39") "0.5853"
$ws.Range("E39").Value = "  -0.83%  "
Set-TextValue $ws.Range("D40") "7.810"
$ws.Range("E40").Value = "  -7.59%  "
$ws.Range("E41").Value = "  -0.31%  "
Set-TextValue $ws.Range("D42") "10.09"
$ws.Range("E42").Value = "  -1.08%  "
Set-TextValue $ws.Range("D43") "2.420"
$ws.Range("E43").Value = "  +3.52%  "
Set-TextValue $ws.Range("D44") "1.286"
$ws.Range("E44").Value = "  +1.92%  "
Set-TextValue $ws.Range("D45") "0.07717"
$ws.Range("E45").Value = "  +2.77%  "
Set-TextValue $ws.Range("D46") "12.25"
$ws.Range("E46").Value = "  +0.30%  "
Set-TextValue $ws.Range("D47") "0.5494"
$ws.Range("E47").Value = "  -1.71%  "
Set-TextValue $ws.Range("D48") "1.918"
$ws.Range("E48").Value = "  -1.04%  "
Set-TextValue $ws.Range("D49") "112.87"
$ws.Range("E49").Value = "  -0.29%  "
Set-TextValue $ws.Range("D50") "0.2951"
$ws.Range("E50").Value = "  -1.16%  "
Set-TextValue $ws.Range("D51") "43.63"
$ws.Range("E51").Value = "  -0.92%  "
